$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.671.36'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.605.50'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.37'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.33%  '

$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.83'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +7.69%  '

$ws.Range("E9").Value = '  +3.08%  '

$ws.Range("E10").Value = '  +1.79%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0908'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.44%  '

$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.600.26'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.07%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.562'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.696.45'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.51%  '

$ws.Range("E16").Value = '  +1.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.38'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.37%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '241.31'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.24'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +8.53%  '

$ws.Range("E20").Value = '  +1.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.998'
$ws.Range("D21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.04'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.49'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.38%  '

$ws.Range("E24").Value = '  +2.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.11'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.62%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.53'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.60%  '

$ws.Range("E27").Value = '  +1.00%  '

$ws.Range("E28").Value = '  +2.09%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("E30").Value = '  +1.69%  '

$ws.Range("E31").Value = '  +0.21%  '

$ws.Range("E32").Value = '  +0.40%  '

$ws.Range("E33").Value = '  +2.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.426.05'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.16%  '

$ws.Range("E35").Value = '  +4.68%  '

$ws.Range("E36").Value = '  +0.27%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.85'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.95%  '

$ws.Range("E38").Value = '  +0.06%  '

$ws.Range("E39").Value = '  +2.81%  '

$ws.Range("E40").Value = '  +3.34%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.97'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.02%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.822'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.52%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0494'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.72%  '

$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '54.39'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '67.83'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.997'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +19.23%  '

$ws.Range("E48").Value = '  +2.97%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.743.89'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.12%  '

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.93'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.54%  '

$ws.Range("B51").Value = 'mCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.12'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.39%  '
